$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "12/19/2025"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = 600.8780000000006
$ws.Range("C19").Value = 0.04118972570139025
$ws.Range("D19").Value = 25
